# The deck's slide order had slide 2 and slide 3 swapped: the slide that
# used to be third (the two-picture "map" slide) now appears second, and
# the slide that used to be second (the black-background single-picture
# slide) now appears third. Reproduce that by moving the slide currently
# at position 3 up to position 2 (PowerPoint's Slide.MoveTo reorders the
# deck just like dragging the slide in the Slides pane).
$p = $ppt.ActivePresentation
$p.Slides.Item(3).MoveTo(2)
